$d = $word.ActiveDocument

# Namespace-qualified OOXML fragments are wrapped in a minimal single-part
# "WordOpenXML" package so Range.InsertXML can parse them; InsertXML replaces
# the contents of the exact range it is called on.
$pkgOpen  = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1. Paragraph 2 ("多云...六一儿童节..."): drop the pPr/rPr/rFonts (the
#        hint=eastAsia formatting that used to be stamped on the paragraph
#        mark itself). Rewrite the whole paragraph without a pPr.
$p2 = $d.Paragraphs(2)
$frag2 = $pkgOpen + `
  '<w:p w:rsidR="00EB1651" w:rsidRDefault="00EB1651" w:rsidP="00EB1651">' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>多云，今天是六一儿童节，又是开心的一天呢</w:t></w:r>' + `
  '</w:p>' + `
  $pkgClose
$null = $p2.Range.InsertXML($frag2)

# --- 2. Paragraph 4 ("中雨...端午节了"): add a pPr/rPr/rFonts hint=eastAsia
#        (the paragraph mark now carries the eastAsia font hint too).
$p4 = $d.Paragraphs(4)
$frag4 = $pkgOpen + `
  '<w:p w:rsidR="00D76E91" w:rsidRDefault="00D76E91" w:rsidP="00D76E91">' + `
    '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>中雨，今天是农历五月初四，明天就是端午节了</w:t></w:r>' + `
  '</w:p>' + `
  $pkgClose
$null = $p4.Range.InsertXML($frag4)

# --- 3. The old bookmark-only paragraph (still paragraph 5) is replaced by
#        two paragraphs: a brand-new "2022年6月3日星期五" paragraph, and a
#        new "中雨...端午节" paragraph that keeps the original bookmark.
$p5 = $d.Paragraphs(5)
$frag5 = $pkgOpen + `
  '<w:p>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2022</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>年</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>6</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>月</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>3</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>日星期五</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p w:rsidR="00D76E91" w:rsidRPr="00D76E91" w:rsidRDefault="00D76E91" w:rsidP="00EB1651">' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>中雨，今天是农历五月初五，中国传统</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>端午节</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>' + `
  $pkgClose
$null = $p5.Range.InsertXML($frag5)

# --- 4. Two brand-new empty paragraphs before the final (already empty)
#        paragraph.
$last = $d.Paragraphs($d.Paragraphs.Count)
$rLast = $last.Range
$rLast.Collapse(1)
$fragEmpty = $pkgOpen + '<w:p/><w:p/><w:p/>' + $pkgClose
$null = $rLast.InsertXML($fragEmpty)
